# 10Th - MB for single stock and added new group
#
# This report tracks, for a single stock (CMCSA), a weekly snapshot of
# each covering analyst firm's most-recent rating action. Every week a
# new column is prepended (after the fixed "UN" status column B) holding
# that week's snapshot for every firm, and any firm whose most-recent
# action falls in one of the new weeks gets its note moved into the
# correct new column. This edit adds three new weekly columns (Jun_26,
# Jun_26, Jun_27 - newest first) and two new analyst rows (Benchmark,
# Evercore ISI).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert 3 new weekly columns right after column B (the fixed "UN"
#     status column). This pushes the old C/D/E (Jun_15/Jun_13/Jun_10)
#     columns to F/G/H, carrying along each cell's value AND style - so
#     the one firm (BidaskClub, row 22) with a real rating note in the
#     oldest-at-the-time column keeps that note+style in its new slot.
$ws.Range("C1:E1").EntireColumn.Insert()

# New columns keep the same 8.0 width as the existing weekly columns.
$ws.Columns.Item(3).ColumnWidth = 8.0
$ws.Columns.Item(4).ColumnWidth = 8.0
$ws.Columns.Item(5).ColumnWidth = 8.0

# --- Header row: newest week goes in B1 (the old B1 "Jun_17" value was
#     just shifted to E1 by the column insert above), and the two
#     brand-new intermediate weeks share the "Jun_26" label in C1/D1.
$ws.Range("B1").Value = "Jun_27"
$ws.Range("C1").Value = "Jun_26"
$ws.Range("D1").Value = "Jun_26"

# --- Body rows: the newly inserted C/D/E cells come back empty from the
#     column insert, so every firm row needs the default "UN" (= no
#     change that week) filled back in, matching the rest of the table.
$lastRow = 27
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Range("C$r").Value = "UN"
    $ws.Range("D$r").Value = "UN"
    $ws.Range("E$r").Value = "UN"
}

# --- Two new analyst-firm rows added at the bottom of the table. Only
#     the fixed firm-name/status columns (A-D) are populated for them,
#     same as a freshly-added firm with no historical weekly data yet.
$ws.Range("A28").Value = "Benchmark"
$ws.Range("B28").Value = "UN"
$ws.Range("C28").Value = "UN"
$ws.Range("D28").Value = "UN"

$ws.Range("A29").Value = "Evercore ISI"
$ws.Range("B29").Value = "UN"
$ws.Range("C29").Value = "UN"
$ws.Range("D29").Value = "UN"
